$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "regFName"
$ws.Range("A7").Value = "regLName"
$ws.Range("B6").Value = "Jordi"
$ws.Range("B7").Value = "D"
$ws.Range("A8").Value = "regUserName"
$ws.Range("B8").Value = "JordiD"
$ws.Range("A9").Value = "regPassword"
$ws.Range("B9").Value = "RamiMaHakesher1#"
$ws.Range("A10").Value = "regWrongPassword"
$ws.Range("B10").Value = "Ya111111"
$ws.Range("A11").Value = "regWrongPasswordMSG"
$ws.Range("B11").Value = "Passwords must have at least one non alphanumeric character, one digit ('0'-'9'), one uppercase ('A'-'Z'), one lowercase ('a'-'z'), one special character and Password must be eight characters or longer."
$ws.Range("A12").Value = "regCaptchaMSG"
$ws.Range("B12").Value = "Please verify reCaptcha to register!"

$fc = $ws.Range("A8:B8").FormatConditions.AddUniqueValues()
$fc.DupeUnique = 1
$fc.Font.Color = 393372
$fc.Interior.Color = 13551615

[void]$ws.Range("B12").Select()
